$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Find.Execute("2024-12-15 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-16 Monday", 2) | Out-Null

# Update table cell values (row-major order matching document order)
$t = $d.Tables.Item(1)
$newValues = @(
  "77-76=",
  "47-12=",
  "50+14=",
  "18+30=",
  "73-68=",
  "61-49=",
  "21-21=",
  "73-72=",
  "8+88=",
  "11+14=",
  "34+20=",
  "61-1=",
  "36+9=",
  "67-6=",
  "66-56=",
  "89-15=",
  "94-90=",
  "45+40=",
  "42-28=",
  "39-31=",
  "96-84=",
  "4+46=",
  "70-36=",
  "26+52=",
  "51-8=",
  "58+38=",
  "83-44=",
  "75+1=",
  "42-41=",
  "5+24=",
  "5+63=",
  "4+42=",
  "33-24=",
  "35-14=",
  "0+50=",
  "49+20=",
  "57-52=",
  "38-6=",
  "54+0=",
  "73-32=",
  "55-1=",
  "51+15=",
  "22+76=",
  "92-61=",
  "36-4=",
  "46+26=",
  "83-82=",
  "39-23=",
  "91-81=",
  "13+58=",
  "48-24=",
  "23+5=",
  "22+1=",
  "17-1=",
  "16+72=",
  "30+64=",
  "44-29=",
  "38+20=",
  "21+5=",
  "19+13=",
  "78-60=",
  "74-9=",
  "20+1=",
  "33+51=",
  "50+41=",
  "67-3=",
  "52+18=",
  "14+68=",
  "13+0=",
  "73-46=",
  "62-47=",
  "84-54=",
  "78+8=",
  "28-7=",
  "45+0=",
  "28+35=",
  "67-18=",
  "52-48=",
  "29+39=",
  "44-14=",
  "61-59=",
  "19+24=",
  "53-17=",
  "38-34=",
  "59+6=",
  "47+50=",
  "69-53=",
  "74-47=",
  "83+15=",
  "67-49=",
  "46+15=",
  "8+40=",
  "31-1=",
  "56-54=",
  "19+47=",
  "61+38=",
  "77-11=",
  "48-4=",
  "55-54=",
  "19+74="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
  for ($c = 1; $c -le 5; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $newValues[$idx]
    $idx++
  }
}

Write-Host "Done. idx=" $idx